$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update the two header text strings (Volume/Number and the report dates)
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/16/2024  Through  9/22/2024"

# ---------------------------------------------------------------------------
# 2) Fix up the cells that flip between "N/A-style text" and real numbers so
#    that they end up with the same number-formatted style as their siblings
#    in the same row (Excel keeps the old style when merely assigning a
#    numeric Value to a text cell, so we first "stamp" the desired style by
#    copying a same-row donor cell that already has it, then overwrite the
#    value).
# ---------------------------------------------------------------------------
$ws.Range("G31").Copy($ws.Range("D31"))   # -> numeric style (#,##0)
$ws.Range("H31").Copy($ws.Range("E31"))   # -> numeric style (#,##0.0)
$ws.Range("C31").Copy($ws.Range("C33"))   # -> text style showing "0"
$ws.Range("G31").Copy($ws.Range("D33"))   # -> numeric style (#,##0)
$ws.Range("H31").Copy($ws.Range("E33"))   # -> numeric style (#,##0.0)

# ---------------------------------------------------------------------------
# 3) Write the updated crime-statistics figures
# ---------------------------------------------------------------------------
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 100
$ws.Range("G14").Value = 12
$ws.Range("H14").Value = -25
$ws.Range("I14").Value = 84
$ws.Range("J14").Value = 102
$ws.Range("K14").Value = -17.647058823529
$ws.Range("L14").Value = -16
$ws.Range("M14").Value = -13.402061855670
$ws.Range("N14").Value = -77.111716621253
$ws.Range("F15").Value = 35
$ws.Range("G15").Value = 25
$ws.Range("H15").Value = 40
$ws.Range("I15").Value = 305
$ws.Range("J15").Value = 283
$ws.Range("K15").Value = 7.773851590106
$ws.Range("L15").Value = 4.810996563573
$ws.Range("M15").Value = 35.555555555555
$ws.Range("N15").Value = -42.015209125475
$ws.Range("C16").Value = 98
$ws.Range("D16").Value = 95
$ws.Range("E16").Value = 3.157894736842
$ws.Range("F16").Value = 405
$ws.Range("G16").Value = 399
$ws.Range("H16").Value = 1.503759398496
$ws.Range("I16").Value = 3628
$ws.Range("J16").Value = 3541
$ws.Range("K16").Value = 2.456933069754
$ws.Range("L16").Value = -2.368137782561
$ws.Range("M16").Value = 15.321042593769
$ws.Range("N16").Value = -68.788713007570
$ws.Range("C17").Value = 164
$ws.Range("D17").Value = 146
$ws.Range("E17").Value = 12.328767123287
$ws.Range("F17").Value = 683
$ws.Range("G17").Value = 698
$ws.Range("H17").Value = -2.148997134670
$ws.Range("I17").Value = 6126
$ws.Range("J17").Value = 6016
$ws.Range("K17").Value = 1.828457446808
$ws.Range("L17").Value = 13.15108976727
$ws.Range("M17").Value = 86.768292682926
$ws.Range("N17").Value = -8.988263259545
$ws.Range("C18").Value = 54
$ws.Range("D18").Value = 55
$ws.Range("E18").Value = -1.818181818181
$ws.Range("F18").Value = 220
$ws.Range("G18").Value = 225
$ws.Range("H18").Value = -2.222222222222
$ws.Range("I18").Value = 2118
$ws.Range("J18").Value = 2178
$ws.Range("K18").Value = -2.754820936639
$ws.Range("L18").Value = -0.981767180925
$ws.Range("M18").Value = -10.557432432432
$ws.Range("N18").Value = -84.514147839438
$ws.Range("C19").Value = 197
$ws.Range("D19").Value = 148
$ws.Range("E19").Value = 33.108108108108
$ws.Range("F19").Value = 771
$ws.Range("G19").Value = 681
$ws.Range("H19").Value = 13.215859030837
$ws.Range("I19").Value = 6754
$ws.Range("J19").Value = 5832
$ws.Range("K19").Value = 15.809327846364
$ws.Range("L19").Value = 16.068052930056
$ws.Range("M19").Value = 99.350649350649
$ws.Range("N19").Value = 22.355072463768
$ws.Range("C20").Value = 95
$ws.Range("D20").Value = 101
$ws.Range("E20").Value = -5.940594059405
$ws.Range("F20").Value = 362
$ws.Range("G20").Value = 403
$ws.Range("H20").Value = -10.173697270471
$ws.Range("I20").Value = 3103
$ws.Range("J20").Value = 3896
$ws.Range("K20").Value = -20.354209445585
$ws.Range("L20").Value = 9.762999646268
$ws.Range("M20").Value = 105.632869449967
$ws.Range("N20").Value = -72.092814101987
$ws.Range("C21").Value = 620
$ws.Range("D21").Value = 552
$ws.Range("E21").Value = 12.318840579710
$ws.Range("F21").Value = 2485
$ws.Range("G21").Value = 2443
$ws.Range("H21").Value = 1.719197707736
$ws.Range("I21").Value = 22118
$ws.Range("J21").Value = 21848
$ws.Range("K21").Value = 1.235811058220
$ws.Range("L21").Value = 8.923470895301
$ws.Range("M21").Value = 57.839149361307
$ws.Range("N21").Value = -55.374868856428
$ws.Range("C22").Value = 4
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -20
$ws.Range("F22").Value = 28
$ws.Range("G22").Value = 26
$ws.Range("H22").Value = 7.692307692307
$ws.Range("I22").Value = 246
$ws.Range("J22").Value = 216
$ws.Range("K22").Value = 13.888888888888
$ws.Range("L22").Value = -8.208955223880
$ws.Range("M22").Value = 8.849557522123
$ws.Range("C23").Value = 31
$ws.Range("D23").Value = 29
$ws.Range("E23").Value = 6.896551724137
$ws.Range("F23").Value = 133
$ws.Range("G23").Value = 138
$ws.Range("H23").Value = -3.623188405797
$ws.Range("I23").Value = 1270
$ws.Range("J23").Value = 1316
$ws.Range("K23").Value = -3.495440729483
$ws.Range("L23").Value = 8.177172061328
$ws.Range("M23").Value = 64.082687338501
$ws.Range("C24").Value = 299
$ws.Range("D24").Value = 409
$ws.Range("E24").Value = -26.894865525672
$ws.Range("F24").Value = 1276
$ws.Range("G24").Value = 1461
$ws.Range("H24").Value = -12.662559890486
$ws.Range("I24").Value = 11814
$ws.Range("J24").Value = 13143
$ws.Range("K24").Value = -10.111846610362
$ws.Range("L24").Value = -12.573077776955
$ws.Range("M24").Value = 27.402135231316
$ws.Range("C25").Value = 108
$ws.Range("D25").Value = 150
$ws.Range("E25").Value = -28
$ws.Range("F25").Value = 470
$ws.Range("G25").Value = 612
$ws.Range("H25").Value = -23.202614379085
$ws.Range("I25").Value = 4684
$ws.Range("J25").Value = 5598
$ws.Range("K25").Value = -16.327259735619
$ws.Range("L25").Value = -31.878999418266
$ws.Range("C26").Value = 215
$ws.Range("D26").Value = 188
$ws.Range("E26").Value = 14.361702127659
$ws.Range("F26").Value = 870
$ws.Range("G26").Value = 795
$ws.Range("H26").Value = 9.433962264150
$ws.Range("I26").Value = 8152
$ws.Range("J26").Value = 7652
$ws.Range("K26").Value = 6.534239414532
$ws.Range("L26").Value = 11.168689485885
$ws.Range("M26").Value = -0.049043648847
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = 15
$ws.Range("E27").Value = -33.333333333333
$ws.Range("I27").Value = 463
$ws.Range("J27").Value = 481
$ws.Range("K27").Value = -3.742203742203
$ws.Range("L27").Value = -8.316831683168
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 32
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 85
$ws.Range("G28").Value = 99
$ws.Range("H28").Value = -14.141414141414
$ws.Range("I28").Value = 873
$ws.Range("J28").Value = 771
$ws.Range("K28").Value = 13.229571984435
$ws.Range("L28").Value = 30.493273542600
$ws.Range("C29").Value = 9
$ws.Range("D29").Value = 7
$ws.Range("E29").Value = 28.571428571428
$ws.Range("F29").Value = 41
$ws.Range("G29").Value = 34
$ws.Range("H29").Value = 20.588235294117
$ws.Range("I29").Value = 323
$ws.Range("J29").Value = 302
$ws.Range("K29").Value = 6.953642384105
$ws.Range("L29").Value = -15.223097112860
$ws.Range("M29").Value = -10.773480662983
$ws.Range("N29").Value = -69.001919385796
$ws.Range("C30").Value = 7
$ws.Range("D30").Value = 6
$ws.Range("E30").Value = 16.666666666666
$ws.Range("F30").Value = 33
$ws.Range("G30").Value = 29
$ws.Range("H30").Value = 13.793103448275
$ws.Range("I30").Value = 254
$ws.Range("J30").Value = 248
$ws.Range("K30").Value = 2.419354838709
$ws.Range("L30").Value = -20.625
$ws.Range("M30").Value = -15.894039735099
$ws.Range("N30").Value = -72.978723404255
$ws.Range("G31").Value = 2
$ws.Range("J31").Value = 16
$ws.Range("K31").Value = 0
$ws.Range("F33").Value = 3
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = -25
$ws.Range("J33").Value = 35
$ws.Range("K33").Value = 2.857142857142
$ws.Range("L33").Value = -8.316831683168

$ws.Range("D31").Value = 1
$ws.Range("E31").Value = -100
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
